$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -11.9692
$ws.Range("B7").Value = 4.793099999999997
$ws.Range("A8").Value = -22.39200000000002
$ws.Range("A10").Value = -21.96250000000001
$ws.Range("D10").Value = -7.940399999999999
$ws.Range("A12").Value = -21.53380000000001
$ws.Range("D12").Value = -7.231400000000001
$ws.Range("D13").Value = -8.96409999999999
$ws.Range("D14").Value = -7.907099999999997
$ws.Range("B15").Value = 5.187199999999994
$ws.Range("A18").Value = -21.80099999999998
$ws.Range("B18").Value = 6.195799999999998
$ws.Range("C18").Value = -11.9167
$ws.Range("C19").Value = -11.5915
$ws.Range("B20").Value = 8.869100000000007
$ws.Range("C27").Value = -13.36489999999999
$ws.Range("B29").Value = 5.0267
$ws.Range("D29").Value = -7.449799999999995
$ws.Range("B30").Value = 5.4946
$ws.Range("B31").Value = 4.8776
$ws.Range("C31").Value = -13.65409999999999
$ws.Range("D32").Value = -9.017499999999998
$ws.Range("D35").Value = -7.750799999999999
$ws.Range("A37").Value = -19.22859999999999
$ws.Range("C38").Value = -13.0476
$ws.Range("B40").Value = 9.514499999999995
$ws.Range("C42").Value = -12.0799
$ws.Range("D43").Value = -8.154300000000005
$ws.Range("C44").Value = -13.22349999999999
$ws.Range("C47").Value = -12.4627
$ws.Range("D48").Value = -7.484199999999995
$ws.Range("D49").Value = -8.248800000000001
$ws.Range("B50").Value = 5.211799999999998
$ws.Range("D50").Value = -8.191699999999992
$ws.Range("A55").Value = -22.5139
$ws.Range("D56").Value = -8.219599999999996
$ws.Range("C58").Value = -12.5623
$ws.Range("C65").Value = -12.4561
$ws.Range("A68").Value = -21.48140000000001
$ws.Range("B68").Value = 4.432600000000001
$ws.Range("D69").Value = -7.242099999999994
$ws.Range("C73").Value = -12.5185
$ws.Range("B76").Value = 6.207899999999998
$ws.Range("A77").Value = -20.20769999999998
$ws.Range("A78").Value = -20.00119999999998
$ws.Range("A81").Value = -21.875
$ws.Range("D81").Value = -7.704099999999998
$ws.Range("A82").Value = -21.9257
$ws.Range("B87").Value = 4.860999999999995
$ws.Range("B88").Value = 4.509499999999998
$ws.Range("C90").Value = -12.9296
$ws.Range("D92").Value = -6.254199999999997
$ws.Range("C94").Value = -10.0879
$ws.Range("C95").Value = -12.0734
$ws.Range("B96").Value = 5.370300000000007
$ws.Range("B98").Value = 6.079400000000001
$ws.Range("B101").Value = 9.534599999999996
$ws.Range("C101").Value = -12.68910000000001
$ws.Range("B102").Value = 8.554100000000004
